# Append the newest daily COVID-19 data point (2020-06-01, published 2020-06-02)
# as a new row at the bottom of the "Tabela1" table on the active sheet,
# growing the table/autofilter/dimension/selection accordingly.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Tabela1")

$lastRow = $tbl.Range.Row + $tbl.Range.Rows.Count - 1   # last row currently used by the table (82)
$newRow  = $lastRow + 1                                  # row to be created (83)

# Insert a new (blank) row at the bottom of the table range, copying the
# formatting (number formats / styles) from the row above it, just like
# typing a new row right under the table in the Excel UI.
$ws.Rows.Item($newRow).Insert(-4121, 0) | Out-Null   # xlShiftDown, xlFormatFromLeftOrAbove

# Grow the table (and its autoFilter) so the new row becomes part of it.
$tbl.Resize($ws.Range($ws.Cells.Item($tbl.Range.Row, $tbl.Range.Column), $ws.Cells.Item($newRow, $tbl.Range.Column + $tbl.Range.Columns.Count - 1))) | Out-Null

# Fill in the new row's data.
$ws.Cells.Item($newRow, 1).Value2  = 43983   # Date            -> 01/06/2020
$ws.Cells.Item($newRow, 2).Value2  = 79698   # Tested (all)
$ws.Cells.Item($newRow, 3).Value2  = 659     # Tested (daily)
$ws.Cells.Item($newRow, 4).Value2  = 1475    # Positive (all)
$ws.Cells.Item($newRow, 5).Value2  = 2       # Positive (daily)
$ws.Cells.Item($newRow, 6).Value2  = 5       # All hospitalized on certain day
$ws.Cells.Item($newRow, 7).Value2  = 1       # All persons in intensive care on certain day
$ws.Cells.Item($newRow, 8).Value2  = 0       # Discharged
$ws.Cells.Item($newRow, 9).Value2  = 109     # Deaths (all)
$ws.Cells.Item($newRow, 10).Value2 = 0       # Deaths (daily)

# Match the author's selection state (row is now the last, fully selected).
$ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 10)).Select() | Out-Null
